# log_pcsmote_x_muestra_heart_D25_R25_Pentropia.xlsx
# "bug guardado corregido y aumento ventana pureza_proporcion"
#
# This re-run of the pcsmote logging pass:
#   - fixed a save bug so that rows whose synthetic-count was previously
#     (incorrectly) logged as 0/blank now carry their real
#     synthetics_from_this_seed / last_delta / last_neighbor_z values
#   - widened the "pureza_proporcion" window, which changed several of the
#     synthetics_from_this_seed (W), last_delta (X) and last_neighbor_z (Y)
#     counts for rows 2-60
#   - refreshed the timestamp (Z) on every surviving row
#   - the window change also means there are no longer extra seed rows
#     61-68 in this run, so those rows are removed entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for W (synthetics_from_this_seed), X (last_delta) and
# Y (last_neighbor_z) - only rows whose numbers actually changed are listed;
# everything else in columns A-Y is left untouched.
$wxyUpdates = @(
    @{ Row = 2;  W = 3;  X = 0.4538824667597043; Y = 30 }
    @{ Row = 3;  W = 1;  X = 0.5616240759128834; Y = 123 }
    @{ Row = 4;  W = 3;  X = 0.4329311706285884; Y = 24 }
    @{ Row = 5;  W = 5;  X = 0.4153959819657586; Y = 6 }
    @{ Row = 8;  W = 6;  X = 0.5079682182603347; Y = 11 }
    @{ Row = 12; W = 4;  X = 0.5614880310328125; Y = 120 }
    @{ Row = 14; W = 1;  X = 0.4062858371373469; Y = 61 }
    @{ Row = 16; W = 4;  X = 0.405083825348819;  Y = 21 }
    @{ Row = 26; W = 3;  X = 0.4705137712668338; Y = 18 }
    @{ Row = 27; W = 3;  X = 0.5541934359909122; Y = 123 }
    @{ Row = 28; W = 2;  X = 0.4969659942717967; Y = 134 }
    @{ Row = 29; W = 4;  X = 0.4636006949943728; Y = 31 }
    @{ Row = 30; W = 9;  X = 0.4641560129943472; Y = 108 }
    @{ Row = 31; W = 7;  X = 0.4978905520555126; Y = 22 }
    @{ Row = 32; W = 9;  X = 0.5290345580818899; Y = 80 }
    @{ Row = 33; W = 7;  X = 0.5381875476204931; Y = 6 }
    @{ Row = 34; W = 5;  X = 0.5966846281789686; Y = 91 }
    @{ Row = 36; W = 9;  X = 0.4081550283109528; Y = 135 }
    @{ Row = 38; W = 13; X = 0.4739308912122809; Y = 11 }
    @{ Row = 41; W = 7;  X = 0.4727259204758588; Y = 118 }
    @{ Row = 44; W = 4;  X = 0.4822074026636463; Y = 16 }
    @{ Row = 45; W = 4;  X = 0.42961738599068;    Y = 116 }
    @{ Row = 46; W = 5;  X = 0.5721461166512687; Y = 16 }
    @{ Row = 49; W = 9;  X = 0.4690142496053366; Y = 115 }
    @{ Row = 50; W = 4;  X = 0.453356202855057;  Y = 52 }
    @{ Row = 51; W = 6;  X = 0.4444215620941461; Y = 17 }
    @{ Row = 52; W = 5;  X = 0.5792182599846987; Y = 44 }
    @{ Row = 57; W = 12; X = 0.5392608545679577; Y = 118 }
    @{ Row = 59; W = 9;  X = 0.450783082786869;  Y = 144 }
    @{ Row = 60; W = 6;  X = 0.4646405864041511; Y = 22 }
)

foreach ($u in $wxyUpdates) {
    $r = $u.Row
    $ws.Range("W$r").Value = $u.W
    $ws.Range("X$r").Value = $u.X
    $ws.Range("Y$r").Value = $u.Y
}

# New timestamp (column Z) for every remaining data row (2-60).
$timestamps = @{
    2  = "2025-10-29T23:40:24.423720"; 3  = "2025-10-29T23:40:24.423720"
    4  = "2025-10-29T23:40:24.423720"; 5  = "2025-10-29T23:40:24.423720"
    6  = "2025-10-29T23:40:24.423720"; 7  = "2025-10-29T23:40:24.423720"
    8  = "2025-10-29T23:40:24.423720"; 9  = "2025-10-29T23:40:24.424719"
    10 = "2025-10-29T23:40:24.424719"; 11 = "2025-10-29T23:40:24.424719"
    12 = "2025-10-29T23:40:24.424719"; 13 = "2025-10-29T23:40:24.424719"
    14 = "2025-10-29T23:40:24.424719"; 15 = "2025-10-29T23:40:24.424719"
    16 = "2025-10-29T23:40:24.424719"; 17 = "2025-10-29T23:40:24.425728"
    18 = "2025-10-29T23:40:24.425728"; 19 = "2025-10-29T23:40:24.425728"
    20 = "2025-10-29T23:40:24.425728"; 21 = "2025-10-29T23:40:24.425728"
    22 = "2025-10-29T23:40:24.425728"; 23 = "2025-10-29T23:40:24.425728"
    24 = "2025-10-29T23:40:24.425728"; 25 = "2025-10-29T23:40:24.426719"
    26 = "2025-10-29T23:40:24.426719"; 27 = "2025-10-29T23:40:24.426719"
    28 = "2025-10-29T23:40:24.426719"; 29 = "2025-10-29T23:40:24.426719"
    30 = "2025-10-29T23:40:24.464051"; 31 = "2025-10-29T23:40:24.464051"
    32 = "2025-10-29T23:40:24.465051"; 33 = "2025-10-29T23:40:24.465051"
    34 = "2025-10-29T23:40:24.465051"; 35 = "2025-10-29T23:40:24.465051"
    36 = "2025-10-29T23:40:24.466054"; 37 = "2025-10-29T23:40:24.466054"
    38 = "2025-10-29T23:40:24.466054"; 39 = "2025-10-29T23:40:24.467052"
    40 = "2025-10-29T23:40:24.467052"; 41 = "2025-10-29T23:40:24.469055"
    42 = "2025-10-29T23:40:24.469055"; 43 = "2025-10-29T23:40:24.511505"
    44 = "2025-10-29T23:40:24.512505"; 45 = "2025-10-29T23:40:24.512505"
    46 = "2025-10-29T23:40:24.512505"; 47 = "2025-10-29T23:40:24.512505"
    48 = "2025-10-29T23:40:24.512505"; 49 = "2025-10-29T23:40:24.512505"
    50 = "2025-10-29T23:40:24.512505"; 51 = "2025-10-29T23:40:24.512505"
    52 = "2025-10-29T23:40:24.512505"; 53 = "2025-10-29T23:40:24.512505"
    54 = "2025-10-29T23:40:24.512505"; 55 = "2025-10-29T23:40:24.513505"
    56 = "2025-10-29T23:40:24.513505"; 57 = "2025-10-29T23:40:24.513505"
    58 = "2025-10-29T23:40:24.513505"; 59 = "2025-10-29T23:40:24.513505"
    60 = "2025-10-29T23:40:24.513505"
}

foreach ($r in 2..60) {
    $ws.Range("Z$r").Value = $timestamps[$r]
}

# The wider pureza_proporcion window means this run produced fewer sampled
# seed rows - rows 61-68 from the previous run no longer exist.
$ws.Rows("61:68").Delete()

# Shrink the conditional-formatting range so it still covers exactly the
# data rows (A2:Z68 -> A2:Z60).
$cond = $ws.Range("A2:Z68").FormatConditions.Item(1)
$cond.ModifyAppliesToRange($ws.Range("A2:Z60"))
